$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs | Gal | Galr3 | ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gal"
$ws.Range("C2").Value = "Galr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03644533333333334
$ws.Range("H2").Value = 0.109336
$ws.Range("I2").Value = 0.005561955322140003
$ws.Range("J2").Value = 0.005561955322140003
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03081733333333334
$ws.Range("N2").Value = 0.09245200000000001
$ws.Range("O2").Value = 0.2575601819736846
$ws.Range("P2").Value = 0.2575601819736846
$ws.Range("Q2").Value = 0.001123147985777778
$ws.Range("R2").Value = 0.010108331872
$ws.Range("S2").Value = 0.001432538224899883
$ws.Range("T2").Value = 0.001432538224899883

# Row 3: ECs | Gal | Galr3 | FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gal"
$ws.Range("C3").Value = "Galr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03644533333333334
$ws.Range("H3").Value = 0.109336
$ws.Range("I3").Value = 0.005561955322140003
$ws.Range("J3").Value = 0.005561955322140003
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.08883366666666666
$ws.Range("N3").Value = 0.266501
$ws.Range("O3").Value = 0.7424398180263154
$ws.Range("P3").Value = 0.7424398180263154
$ws.Range("Q3").Value = 0.003237572592888889
$ws.Range("R3").Value = 0.029138153336
$ws.Range("S3").Value = 0.00412941709724012
$ws.Range("T3").Value = 0.00412941709724012

# Row 4: MuSCs | Gal | Galr3 | ECs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Gal"
$ws.Range("C4").Value = "Galr3"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.516166333333334
$ws.Range("H4").Value = 19.548499
$ws.Range("I4").Value = 0.99443804467786
$ws.Range("J4").Value = 0.99443804467786
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03081733333333334
$ws.Range("N4").Value = 0.09245200000000001
$ws.Range("O4").Value = 0.2575601819736846
$ws.Range("P4").Value = 0.2575601819736846
$ws.Range("Q4").Value = 0.2008108699497778
$ws.Range("R4").Value = 1.807297829548
$ws.Range("S4").Value = 0.2561276437487848
$ws.Range("T4").Value = 0.2561276437487848

# Row 5: MuSCs | Gal | Galr3 | FAPs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gal"
$ws.Range("C5").Value = "Galr3"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.516166333333334
$ws.Range("H5").Value = 19.548499
$ws.Range("I5").Value = 0.99443804467786
$ws.Range("J5").Value = 0.99443804467786
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.08883366666666666
$ws.Range("N5").Value = 0.266501
$ws.Range("O5").Value = 0.7424398180263154
$ws.Range("P5").Value = 0.7424398180263154
$ws.Range("Q5").Value = 0.5788549479998889
$ws.Range("R5").Value = 5.209694531998999
$ws.Range("S5").Value = 0.7383104009290753
$ws.Range("T5").Value = 0.7383104009290753
